# Updated cryptos list on Mon Jul 22 17:41:56 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: D/E columns hold text (price strings w/ thousands separators, and
# padded percentage strings) rather than numbers. A leading "'" forces
# Excel to store the assigned value as literal text instead of coercing
# number-looking strings (e.g. "593.92", "1.00") into numeric values.
$updates = @(
    @{ Row = 2;  D = "67.376.08";  E = "  +0.01%  " },
    @{ Row = 3;  D = "3.469.53";   E = "  -1.03%  " },
    @{ Row = 4;  D = $null;        E = "  -0.02%  " },
    @{ Row = 5;  D = "593.92";     E = "  -1.00%  " },
    @{ Row = 6;  D = "181.01";     E = "  +2.90%  " },
    @{ Row = 7;  D = "0.612";      E = "  +3.89%  " },
    @{ Row = 8;  D = $null;        E = "  +0.02%  " },
    @{ Row = 9;  D = "3.465.95";   E = "  -1.07%  " },
    @{ Row = 10; D = "0.140";      E = "  +6.51%  " },
    @{ Row = 11; D = "6.96";       E = "  -2.55%  " },
    @{ Row = 12; D = $null;        E = "  -0.41%  " },
    @{ Row = 13; D = "4.062.48";   E = "  -1.24%  " },
    @{ Row = 14; D = $null;        E = "  +2.79%  " },
    @{ Row = 15; D = $null;        E = "  -0.53%  " },
    @{ Row = 16; D = "67.365.98";  E = "  +0.00%  " },
    @{ Row = 17; D = $null;        E = "  -1.13%  " },
    @{ Row = 18; D = "3.466.73";   E = "  -1.34%  " },
    @{ Row = 19; D = $null;        E = "  -1.69%  " },
    @{ Row = 20; D = "14.14";      E = "  -3.28%  " },
    @{ Row = 21; D = "394.95";     E = "  +0.52%  " },
    @{ Row = 22; D = "7.94";       E = "  -0.72%  " },
    @{ Row = 23; D = "5.79";       E = "  +1.34%  " },
    @{ Row = 24; D = $null;        E = "  +0.02%  " },
    @{ Row = 25; D = $null;        E = "  -0.31%  " },
    @{ Row = 26; D = "71.78";      E = $null },
    @{ Row = 27; D = $null;        E = "  -0.30%  " },
    @{ Row = 28; D = "10.37";      E = "  +0.56%  " },
    @{ Row = 29; D = $null;        E = "  -2.63%  " },
    @{ Row = 30; D = "1.00";       E = "  +0.49%  " },
    @{ Row = 31; D = "6.13";       E = "  +0.07%  " },
    @{ Row = 32; D = "1.41";       E = "  -1.45%  " },
    @{ Row = 33; D = $null;        E = "  -0.83%  " },
    @{ Row = 34; D = "23.54";      E = "  -0.62%  " },
    @{ Row = 35; D = "7.33";       E = "  -0.90%  " },
    @{ Row = 36; D = $null;        E = "  -0.01%  " },
    @{ Row = 37; D = "1.59";       E = "  -3.53%  " },
    @{ Row = 38; D = "160.97";     E = "  -1.64%  " },
    @{ Row = 39; D = "0.880";      E = "  +0.12%  " },
    @{ Row = 40; D = "2.84";       E = "  +11.52%  " },
    @{ Row = 41; D = $null;        E = "  -3.75%  " },
    @{ Row = 42; D = "4.67";       E = "  -0.04%  " },
    @{ Row = 43; D = "6.73";       E = "  -3.84%  " },
    @{ Row = 44; D = "26.16";      E = "  -1.76%  " },
    @{ Row = 45; D = "0.0721";     E = "  -1.45%  " },
    @{ Row = 46; D = "2.760.16";   E = "  -1.82%  " },
    @{ Row = 47; D = "26.36";      E = "  -3.03%  " },
    @{ Row = 48; D = "41.45";      E = "  -2.62%  " },
    @{ Row = 49; D = "0.0299";     E = "  -0.61%  " },
    @{ Row = 50; D = "326.02";     E = "  -3.65%  " },
    @{ Row = 51; D = "1.05";       E = "  -3.18%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $ws.Range("D$r").Value = "'" + $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E$r").Value = "'" + $u.E
    }
}
